$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each text-cell update, forcing Text type (matching the source
# inlineStr cells) instead of letting Excel auto-coerce numeric-looking
# strings ("1.00", "0.0686", ...) into Number cells, then clear the
# temporary "@" number-format override so the cell style index reverts
# to the original (unstyled/General) state.
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value2 = "67.836.39"
$cell.ClearFormats()

$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value2 = "  +1.27%  "
$cell.ClearFormats()

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value2 = "3.250.64"
$cell.ClearFormats()

$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value2 = "  +0.09%  "
$cell.ClearFormats()

$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value2 = "  +0.03%  "
$cell.ClearFormats()

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value2 = "583.29"
$cell.ClearFormats()

$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value2 = "  +0.82%  "
$cell.ClearFormats()

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value2 = "183.13"
$cell.ClearFormats()

$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value2 = "  +3.62%  "
$cell.ClearFormats()

$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value2 = "  +0.08%  "
$cell.ClearFormats()

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value2 = "0.599"
$cell.ClearFormats()

$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value2 = "  -1.26%  "
$cell.ClearFormats()

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value2 = "0.134"
$cell.ClearFormats()

$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value2 = "  +4.01%  "
$cell.ClearFormats()

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value2 = "6.68"
$cell.ClearFormats()

$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value2 = "  -0.55%  "
$cell.ClearFormats()

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value2 = "0.415"
$cell.ClearFormats()

$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value2 = "  +1.81%  "
$cell.ClearFormats()

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value2 = "3.812.94"
$cell.ClearFormats()

$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value2 = "  +0.02%  "
$cell.ClearFormats()

$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value2 = "  +0.41%  "
$cell.ClearFormats()

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value2 = "28.55"
$cell.ClearFormats()

$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value2 = "  +2.51%  "
$cell.ClearFormats()

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value2 = "67.832.62"
$cell.ClearFormats()

$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value2 = "  +1.36%  "
$cell.ClearFormats()

$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value2 = "  +2.33%  "
$cell.ClearFormats()

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value2 = "3.262.39"
$cell.ClearFormats()

$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value2 = "  +0.43%  "
$cell.ClearFormats()

$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value2 = "  +0.44%  "
$cell.ClearFormats()

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value2 = "13.56"
$cell.ClearFormats()

$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value2 = "  +1.75%  "
$cell.ClearFormats()

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value2 = "380.69"
$cell.ClearFormats()

$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value2 = "  +3.17%  "
$cell.ClearFormats()

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value2 = "7.65"
$cell.ClearFormats()

$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value2 = "  +2.28%  "
$cell.ClearFormats()

$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value2 = "  +0.11%  "
$cell.ClearFormats()

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value2 = "71.22"
$cell.ClearFormats()

$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value2 = "  +1.08%  "
$cell.ClearFormats()

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value2 = "0.512"
$cell.ClearFormats()

$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value2 = "  +1.13%  "
$cell.ClearFormats()

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value2 = "0.0000119"
$cell.ClearFormats()

$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value2 = "  +0.85%  "
$cell.ClearFormats()

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value2 = "9.82"
$cell.ClearFormats()

$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value2 = "  +0.49%  "
$cell.ClearFormats()

$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value2 = "  +1.49%  "
$cell.ClearFormats()

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value2 = "1.00"
$cell.ClearFormats()

$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value2 = "  -0.27%  "
$cell.ClearFormats()

$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value2 = "  -0.15%  "
$cell.ClearFormats()

$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value2 = "  +0.63%  "
$cell.ClearFormats()

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value2 = "22.87"
$cell.ClearFormats()

$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value2 = "  +1.79%  "
$cell.ClearFormats()

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value2 = "7.15"
$cell.ClearFormats()

$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value2 = "  +5.75%  "
$cell.ClearFormats()

$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value2 = "  +0.01%  "
$cell.ClearFormats()

$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value2 = "  +2.79%  "
$cell.ClearFormats()

$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value2 = "  +2.67%  "
$cell.ClearFormats()

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value2 = "161.27"
$cell.ClearFormats()

$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value2 = "  -6.55%  "
$cell.ClearFormats()

$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value2 = "  -1.78%  "
$cell.ClearFormats()

$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value2 = "  -0.21%  "
$cell.ClearFormats()

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value2 = "26.50"
$cell.ClearFormats()

$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value2 = "  -1.23%  "
$cell.ClearFormats()

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value2 = "6.70"
$cell.ClearFormats()

$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value2 = "  +4.69%  "
$cell.ClearFormats()

$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value2 = "  +6.84%  "
$cell.ClearFormats()

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value2 = "2.58"
$cell.ClearFormats()

$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value2 = "  -0.57%  "
$cell.ClearFormats()

$cell = $ws.Range("B43")
$cell.NumberFormat = "@"
$cell.Value2 = "OKB"
$cell.ClearFormats()

$cell = $ws.Range("C43")
$cell.NumberFormat = "@"
$cell.Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$cell.ClearFormats()

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value2 = "41.26"
$cell.ClearFormats()

$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value2 = "  +2.21%  "
$cell.ClearFormats()

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value2 = "346.52"
$cell.ClearFormats()

$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value2 = "  +1.98%  "
$cell.ClearFormats()

$cell = $ws.Range("B45")
$cell.NumberFormat = "@"
$cell.Value2 = "InjectiveProtocol"
$cell.ClearFormats()

$cell = $ws.Range("C45")
$cell.NumberFormat = "@"
$cell.Value2 = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$cell.ClearFormats()

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value2 = "25.40"
$cell.ClearFormats()

$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value2 = "  +3.67%  "
$cell.ClearFormats()

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value2 = "0.0686"
$cell.ClearFormats()

$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value2 = "  +1.97%  "
$cell.ClearFormats()

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value2 = "2.629.62"
$cell.ClearFormats()

$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value2 = "  -3.16%  "
$cell.ClearFormats()

$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value2 = "  +1.43%  "
$cell.ClearFormats()

$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value2 = "  -0.88%  "
$cell.ClearFormats()

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value2 = "0.992"
$cell.ClearFormats()

$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value2 = "  +0.92%  "
$cell.ClearFormats()

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value2 = "6.20"
$cell.ClearFormats()

$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value2 = "  +3.09%  "
$cell.ClearFormats()

